$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Row 2 - food / 100 / date
$ws.Range("A2").Value = "food"
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 45867.229537037034
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Row 3 - rent / 500 / date
$ws.Range("A3").Value = "rent"
$ws.Range("B3").Value = 500
$ws.Range("C3").Value = 45866.229537037034

# Copy C2's date style onto C3 so both share the same cell style index
# instead of allocating a second, duplicate cellXf entry.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)  # xlPasteFormats
